$d = $word.ActiveDocument

# Paragraphs for "14º slide" ... "28º slide" are paragraphs 35-49 (1-based,
# via $d.Paragraphs). Apply yellow highlighting to every run in each of
# these paragraphs.
for ($i = 35; $i -le 49; $i++) {
    $d.Paragraphs.Item($i).Range.HighlightColorIndex = 7
}

# The "19º slide" paragraph (item 40) additionally got its paragraph-mark
# run properties (w:pPr/w:rPr) highlighted - this happens when the
# highlighted selection spans into the following paragraph. Reproduce this
# exactly by replacing the paragraph's XML with an equivalent version that
# also carries the highlighted paragraph mark.
$p19 = $d.Paragraphs.Item(40)
$r19 = $p19.Range.Duplicate
$xml19 = @'
<w:p w14:paraId="78DA48A8" w14:textId="0D46348F" w:rsidR="00E70F19" w:rsidRDefault="00E70F19" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">19º slide – Nº de artigos publicados...2000-2018 – gráfico 63 da </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>pg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> 149</w:t></w:r></w:p>
'@
$r19.InsertXML($xml19)

Write-Host "Applied yellow highlight to slides 14-28 (paragraphs 35-49)."
